$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.249.95'
$ws.Cells.Item(2, 5).Value = '  -0.48%  '
$ws.Cells.Item(3, 4).Value = '2.351.67'
$ws.Cells.Item(3, 5).Value = '  +4.73%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '234.16'
$ws.Cells.Item(5, 5).Value = '  +1.59%  '
$ws.Cells.Item(6, 5).Value = '  +1.40%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '71.42'
$ws.Cells.Item(7, 5).Value = '  +12.06%  '
$ws.Cells.Item(8, 5).Value = '  +0.05%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.499'
$ws.Cells.Item(9, 5).Value = '  +12.61%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0972'
$ws.Cells.Item(10, 5).Value = '  +1.50%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '27.37'
$ws.Cells.Item(11, 5).Value = '  -1.11%  '
$ws.Cells.Item(12, 5).Value = '  +2.12%  '
$ws.Cells.Item(13, 4).Value = '2.707.92'
$ws.Cells.Item(13, 5).Value = '  +5.03%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '16.14'
$ws.Cells.Item(14, 5).Value = '  +4.49%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '6.32'
$ws.Cells.Item(15, 5).Value = '  +3.91%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.862'
$ws.Cells.Item(16, 5).Value = '  +4.26%  '
$ws.Cells.Item(17, 4).Value = '2.352.73'
$ws.Cells.Item(17, 5).Value = '  +5.04%  '
$ws.Cells.Item(18, 4).Value = '43.231.48'
$ws.Cells.Item(18, 5).Value = '  -0.47%  '
$ws.Cells.Item(19, 5).Value = '  +3.61%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '6.34'
$ws.Cells.Item(20, 5).Value = '  +4.10%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '74.54'
$ws.Cells.Item(21, 5).Value = '  +2.28%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '250.52'
$ws.Cells.Item(22, 5).Value = '  +1.78%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '3.81'
$ws.Cells.Item(23, 5).Value = '  +3.26%  '
$ws.Cells.Item(24, 5).Value = '  +0.07%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.45'
$ws.Cells.Item(25, 5).Value = '  +1.57%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.25'
$ws.Cells.Item(26, 5).Value = '  -1.17%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.03'
$ws.Cells.Item(27, 5).Value = '  +2.90%  '
$ws.Cells.Item(28, 5).Value = '  +4.18%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '172.47'
$ws.Cells.Item(29, 5).Value = '  -0.12%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.56'
$ws.Cells.Item(30, 5).Value = '  +9.67%  '
$ws.Cells.Item(31, 5).Value = '  +0.80%  '
$ws.Cells.Item(32, 5).Value = '  +2.36%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.98'
$ws.Cells.Item(33, 5).Value = '  +1.98%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0690'
$ws.Cells.Item(34, 5).Value = '  +2.38%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.05'
$ws.Cells.Item(35, 5).Value = '  +3.54%  '
$ws.Cells.Item(37, 2).Value = 'THORChain'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '6.53'
$ws.Cells.Item(37, 5).Value = '  +3.72%  '
$ws.Cells.Item(38, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '2.43'
$ws.Cells.Item(38, 5).Value = '  +6.83%  '
$ws.Cells.Item(39, 5).Value = '  +1.87%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '18.97'
$ws.Cells.Item(40, 5).Value = '  +11.56%  '
$ws.Cells.Item(41, 2).Value = 'BinanceUSD'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.00'
$ws.Cells.Item(41, 5).Value = '  -0.16%  '
$ws.Cells.Item(42, 2).Value = 'FraxShare'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '8.92'
$ws.Cells.Item(42, 5).Value = '  +3.39%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '4.50'
$ws.Cells.Item(43, 5).Value = '  +0.16%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '99.13'
$ws.Cells.Item(44, 5).Value = '  +2.84%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.15'
$ws.Cells.Item(45, 5).Value = '  +9.39%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0957'
$ws.Cells.Item(46, 5).Value = '  +2.12%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.20'
$ws.Cells.Item(47, 5).Value = '  +2.29%  '
$ws.Cells.Item(48, 4).Value = '1.442.43'
$ws.Cells.Item(48, 5).Value = '  -0.52%  '
$ws.Cells.Item(49, 4).Value = '2.578.62'
$ws.Cells.Item(49, 5).Value = '  +5.10%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '2.77'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.000201'
$ws.Cells.Item(51, 5).Value = '  -3.75%  '
